$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row content swaps / permutations (existing rows) ---
# row 14
$ws.Range("B14").Value2 = 6772177
$ws.Range("F14").Value2 = "Aguilas Doradas"
$ws.Range("G14").Value2 = "Alianza Petrolera"
$ws.Range("H14").Value2 = 1
$ws.Range("I14").Value2 = 1
$ws.Range("J14").Value2 = "D"
$ws.Range("K14").Value2 = 2.15
$ws.Range("L14").Value2 = 3.3
$ws.Range("M14").Value2 = 3.5
$ws.Range("N14").Value2 = 2.2
$ws.Range("O14").Value2 = 3.5
$ws.Range("P14").Value2 = 3.2
$ws.Range("Q14").Value2 = -0.25
$ws.Range("R14").Value2 = 1.9
$ws.Range("S14").Value2 = 1.9
$ws.Range("T14").Value2 = 2.75
$ws.Range("U14").Value2 = 1.95
$ws.Range("V14").Value2 = 1.85
$ws.Range("W14").Value2 = -1
$ws.Range("X14").Value2 = 2.5
$ws.Range("Y14").Value2 = -1
$ws.Range("Z14").Value2 = -0.5
$ws.Range("AA14").Value2 = 0.45
$ws.Range("AB14").Value2 = -1
$ws.Range("AC14").Value2 = 0.8500000000000001

# row 15
$ws.Range("B15").Value2 = 6772175
$ws.Range("F15").Value2 = "Atletico Nacional Medellin"
$ws.Range("G15").Value2 = "Deportivo Pasto"
$ws.Range("H15").Value2 = 3
$ws.Range("I15").Value2 = 2
$ws.Range("J15").Value2 = "H"
$ws.Range("K15").Value2 = 1.666
$ws.Range("L15").Value2 = 3.75
$ws.Range("M15").Value2 = 4.5
$ws.Range("N15").Value2 = 1.8
$ws.Range("O15").Value2 = 3.6
$ws.Range("P15").Value2 = 5
$ws.Range("Q15").Value2 = -0.75
$ws.Range("R15").Value2 = 2
$ws.Range("S15").Value2 = 1.85
$ws.Range("T15").Value2 = 2.25
$ws.Range("U15").Value2 = 1.85
$ws.Range("V15").Value2 = 2
$ws.Range("W15").Value2 = 0.8
$ws.Range("X15").Value2 = -1
$ws.Range("Y15").Value2 = -1
$ws.Range("Z15").Value2 = 0.5
$ws.Range("AA15").Value2 = -0.5
$ws.Range("AB15").Value2 = 0.8500000000000001
$ws.Range("AC15").Value2 = -1

# row 208
$ws.Range("B208").Value2 = 7404217
$ws.Range("F208").Value2 = "Alianza Petrolera"
$ws.Range("G208").Value2 = "Deportivo Pereira"
$ws.Range("H208").Value2 = 2
$ws.Range("I208").Value2 = 1
$ws.Range("J208").Value2 = "H"
$ws.Range("K208").Value2 = 1.95
$ws.Range("L208").Value2 = 3.2
$ws.Range("M208").Value2 = 3.75
$ws.Range("N208").Value2 = 1.95
$ws.Range("O208").Value2 = 3.2
$ws.Range("P208").Value2 = 4.75
$ws.Range("Q208").Value2 = -0.5
$ws.Range("R208").Value2 = 1.925
$ws.Range("S208").Value2 = 1.875
$ws.Range("T208").Value2 = 2
$ws.Range("U208").Value2 = 1.825
$ws.Range("V208").Value2 = 1.975
$ws.Range("W208").Value2 = 0.95
$ws.Range("X208").Value2 = -1
$ws.Range("Y208").Value2 = -1
$ws.Range("Z208").Value2 = 0.925
$ws.Range("AA208").Value2 = -1
$ws.Range("AB208").Value2 = 0.825
$ws.Range("AC208").Value2 = -1

# row 209
$ws.Range("B209").Value2 = 7404214
$ws.Range("F209").Value2 = "Boyaca Chico"
$ws.Range("G209").Value2 = "Deportivo Cali"
$ws.Range("H209").Value2 = 1
$ws.Range("I209").Value2 = 1
$ws.Range("J209").Value2 = "D"
$ws.Range("K209").Value2 = 3.2
$ws.Range("L209").Value2 = 3.1
$ws.Range("M209").Value2 = 2.2
$ws.Range("N209").Value2 = 3.6
$ws.Range("O209").Value2 = 3
$ws.Range("P209").Value2 = 2.25
$ws.Range("Q209").Value2 = 0.25
$ws.Range("R209").Value2 = 1.95
$ws.Range("S209").Value2 = 1.9
$ws.Range("T209").Value2 = 2.25
$ws.Range("U209").Value2 = 1.875
$ws.Range("V209").Value2 = 1.975
$ws.Range("W209").Value2 = -1
$ws.Range("X209").Value2 = 2
$ws.Range("Y209").Value2 = -1
$ws.Range("Z209").Value2 = 0.475
$ws.Range("AA209").Value2 = -0.5
$ws.Range("AB209").Value2 = -0.5
$ws.Range("AC209").Value2 = 0.4875

# row 211
$ws.Range("B211").Value2 = 7404212
$ws.Range("F211").Value2 = "Envigado FC"
$ws.Range("G211").Value2 = "Deportivo Pasto"
$ws.Range("H211").Value2 = 1
$ws.Range("I211").Value2 = 1
$ws.Range("J211").Value2 = "D"
$ws.Range("K211").Value2 = 2.6
$ws.Range("L211").Value2 = 2.875
$ws.Range("M211").Value2 = 2.8
$ws.Range("N211").Value2 = 2.8
$ws.Range("O211").Value2 = 3.2
$ws.Range("P211").Value2 = 2.625
$ws.Range("Q211").Value2 = 0
$ws.Range("R211").Value2 = 1.975
$ws.Range("S211").Value2 = 1.875
$ws.Range("T211").Value2 = 2.5
$ws.Range("U211").Value2 = 2.025
$ws.Range("V211").Value2 = 1.825
$ws.Range("W211").Value2 = -1
$ws.Range("X211").Value2 = 2.2
$ws.Range("Y211").Value2 = -1
$ws.Range("Z211").Value2 = 0
$ws.Range("AA211").Value2 = -0
$ws.Range("AB211").Value2 = -1
$ws.Range("AC211").Value2 = 0.825

# row 212
$ws.Range("B212").Value2 = 7404218
$ws.Range("F212").Value2 = "Junior"
$ws.Range("G212").Value2 = "Atletico Huila"
$ws.Range("H212").Value2 = 2
$ws.Range("I212").Value2 = 0
$ws.Range("J212").Value2 = "H"
$ws.Range("K212").Value2 = 1.363
$ws.Range("L212").Value2 = 4.5
$ws.Range("M212").Value2 = 7
$ws.Range("N212").Value2 = 1.3
$ws.Range("O212").Value2 = 5
$ws.Range("P212").Value2 = 12
$ws.Range("Q212").Value2 = -1.5
$ws.Range("R212").Value2 = 1.9
$ws.Range("S212").Value2 = 1.95
$ws.Range("T212").Value2 = 2.75
$ws.Range("U212").Value2 = 2.025
$ws.Range("V212").Value2 = 1.825
$ws.Range("W212").Value2 = 0.3
$ws.Range("X212").Value2 = -1
$ws.Range("Y212").Value2 = -1
$ws.Range("Z212").Value2 = 0.8999999999999999
$ws.Range("AA212").Value2 = -1
$ws.Range("AB212").Value2 = -1
$ws.Range("AC212").Value2 = 0.825

# row 213
$ws.Range("B213").Value2 = 7404522
$ws.Range("F213").Value2 = "La Equidad"
$ws.Range("G213").Value2 = "Millonarios"
$ws.Range("H213").Value2 = 2
$ws.Range("I213").Value2 = 1
$ws.Range("J213").Value2 = "H"
$ws.Range("K213").Value2 = 2.4
$ws.Range("L213").Value2 = 3.1
$ws.Range("M213").Value2 = 2.875
$ws.Range("N213").Value2 = 2.1
$ws.Range("O213").Value2 = 3.1
$ws.Range("P213").Value2 = 3.8
$ws.Range("Q213").Value2 = -0.25
$ws.Range("R213").Value2 = 1.75
$ws.Range("S213").Value2 = 2.05
$ws.Range("T213").Value2 = 2
$ws.Range("U213").Value2 = 1.85
$ws.Range("V213").Value2 = 1.95
$ws.Range("W213").Value2 = 1.1
$ws.Range("X213").Value2 = -1
$ws.Range("Y213").Value2 = -1
$ws.Range("Z213").Value2 = 0.75
$ws.Range("AA213").Value2 = -1
$ws.Range("AB213").Value2 = 0.8500000000000001
$ws.Range("AC213").Value2 = -1

# row 214
$ws.Range("B214").Value2 = 7404260
$ws.Range("F214").Value2 = "Atletico Nacional Medellin"
$ws.Range("G214").Value2 = "Deportes Tolima"
$ws.Range("H214").Value2 = 2
$ws.Range("I214").Value2 = 3
$ws.Range("J214").Value2 = "A"
$ws.Range("K214").Value2 = 2
$ws.Range("L214").Value2 = 3.25
$ws.Range("M214").Value2 = 3.5
$ws.Range("N214").Value2 = 1.75
$ws.Range("O214").Value2 = 3.6
$ws.Range("P214").Value2 = 4.75
$ws.Range("Q214").Value2 = -0.75
$ws.Range("R214").Value2 = 2
$ws.Range("S214").Value2 = 1.8
$ws.Range("T214").Value2 = 2.5
$ws.Range("U214").Value2 = 2
$ws.Range("V214").Value2 = 1.8
$ws.Range("W214").Value2 = -1
$ws.Range("X214").Value2 = -1
$ws.Range("Y214").Value2 = 3.75
$ws.Range("Z214").Value2 = -1
$ws.Range("AA214").Value2 = 0.8
$ws.Range("AB214").Value2 = 1
$ws.Range("AC214").Value2 = -1

# row 215
$ws.Range("B215").Value2 = 7404219
$ws.Range("F215").Value2 = "Union Magdalena"
$ws.Range("G215").Value2 = "Independiente Medellin"
$ws.Range("H215").Value2 = 0
$ws.Range("I215").Value2 = 4
$ws.Range("J215").Value2 = "A"
$ws.Range("K215").Value2 = 3
$ws.Range("L215").Value2 = 3.1
$ws.Range("M215").Value2 = 2.3
$ws.Range("N215").Value2 = 3.6
$ws.Range("O215").Value2 = 3.4
$ws.Range("P215").Value2 = 2.1
$ws.Range("Q215").Value2 = 0.25
$ws.Range("R215").Value2 = 2.025
$ws.Range("S215").Value2 = 1.775
$ws.Range("T215").Value2 = 2.5
$ws.Range("U215").Value2 = 1.85
$ws.Range("V215").Value2 = 1.95
$ws.Range("W215").Value2 = -1
$ws.Range("X215").Value2 = -1
$ws.Range("Y215").Value2 = 1.1
$ws.Range("Z215").Value2 = -1
$ws.Range("AA215").Value2 = 0.7749999999999999
$ws.Range("AB215").Value2 = 0.8500000000000001
$ws.Range("AC215").Value2 = -1

# row 216
$ws.Range("B216").Value2 = 7404215
$ws.Range("F216").Value2 = "America de Cali"
$ws.Range("G216").Value2 = "Atletico Bucaramanga"
$ws.Range("H216").Value2 = 1
$ws.Range("I216").Value2 = 2
$ws.Range("J216").Value2 = "A"
$ws.Range("K216").Value2 = 1.444
$ws.Range("L216").Value2 = 4.5
$ws.Range("M216").Value2 = 6
$ws.Range("N216").Value2 = 1.363
$ws.Range("O216").Value2 = 5
$ws.Range("P216").Value2 = 7.5
$ws.Range("Q216").Value2 = -1.25
$ws.Range("R216").Value2 = 1.775
$ws.Range("S216").Value2 = 2.025
$ws.Range("T216").Value2 = 3
$ws.Range("U216").Value2 = 1.925
$ws.Range("V216").Value2 = 1.875
$ws.Range("W216").Value2 = -1
$ws.Range("X216").Value2 = -1
$ws.Range("Y216").Value2 = 6.5
$ws.Range("Z216").Value2 = -1
$ws.Range("AA216").Value2 = 1.025
$ws.Range("AB216").Value2 = 0
$ws.Range("AC216").Value2 = -0

# row 217
$ws.Range("B217").Value2 = 7404213
$ws.Range("F217").Value2 = "Jaguares de Cordoba"
$ws.Range("G217").Value2 = "Aguilas Doradas"
$ws.Range("H217").Value2 = 0
$ws.Range("I217").Value2 = 1
$ws.Range("J217").Value2 = "A"
$ws.Range("K217").Value2 = 3.25
$ws.Range("L217").Value2 = 3.1
$ws.Range("M217").Value2 = 2.2
$ws.Range("N217").Value2 = 3.6
$ws.Range("O217").Value2 = 3.2
$ws.Range("P217").Value2 = 2.15
$ws.Range("Q217").Value2 = 0.25
$ws.Range("R217").Value2 = 1.975
$ws.Range("S217").Value2 = 1.825
$ws.Range("T217").Value2 = 2
$ws.Range("U217").Value2 = 1.75
$ws.Range("V217").Value2 = 2.05
$ws.Range("W217").Value2 = -1
$ws.Range("X217").Value2 = -1
$ws.Range("Y217").Value2 = 1.15
$ws.Range("Z217").Value2 = -1
$ws.Range("AA217").Value2 = 0.825
$ws.Range("AB217").Value2 = -1
$ws.Range("AC217").Value2 = 1.05

# row 238
$ws.Range("B238").Value2 = 7528604
$ws.Range("F238").Value2 = "Aguilas Doradas"
$ws.Range("G238").Value2 = "Deportivo Cali"
$ws.Range("H238").Value2 = 3
$ws.Range("I238").Value2 = 1
$ws.Range("J238").Value2 = "H"
$ws.Range("K238").Value2 = 1.666
$ws.Range("L238").Value2 = 3.75
$ws.Range("M238").Value2 = 5
$ws.Range("N238").Value2 = 1.363
$ws.Range("O238").Value2 = 5
$ws.Range("P238").Value2 = 9
$ws.Range("Q238").Value2 = -1.25
$ws.Range("R238").Value2 = 1.825
$ws.Range("S238").Value2 = 1.975
$ws.Range("T238").Value2 = 2.75
$ws.Range("U238").Value2 = 1.9
$ws.Range("V238").Value2 = 1.9
$ws.Range("W238").Value2 = 0.363
$ws.Range("X238").Value2 = -1
$ws.Range("Y238").Value2 = -1
$ws.Range("Z238").Value2 = 0.825
$ws.Range("AA238").Value2 = -1
$ws.Range("AB238").Value2 = 0.8999999999999999
$ws.Range("AC238").Value2 = -1

# row 239
$ws.Range("B239").Value2 = 7528136
$ws.Range("F239").Value2 = "Millonarios"
$ws.Range("G239").Value2 = "Atletico Nacional Medellin"
$ws.Range("H239").Value2 = 0
$ws.Range("I239").Value2 = 1
$ws.Range("J239").Value2 = "A"
$ws.Range("K239").Value2 = 1.85
$ws.Range("L239").Value2 = 3.3
$ws.Range("M239").Value2 = 4.5
$ws.Range("N239").Value2 = 1.85
$ws.Range("O239").Value2 = 3.5
$ws.Range("P239").Value2 = 4.2
$ws.Range("Q239").Value2 = -0.5
$ws.Range("R239").Value2 = 1.875
$ws.Range("S239").Value2 = 1.975
$ws.Range("T239").Value2 = 2.5
$ws.Range("U239").Value2 = 2.05
$ws.Range("V239").Value2 = 1.8
$ws.Range("W239").Value2 = -1
$ws.Range("X239").Value2 = -1
$ws.Range("Y239").Value2 = 3.2
$ws.Range("Z239").Value2 = -1
$ws.Range("AA239").Value2 = 0.9750000000000001
$ws.Range("AB239").Value2 = -1
$ws.Range("AC239").Value2 = 0.8

# row 240
$ws.Range("B240").Value2 = 7528603
$ws.Range("F240").Value2 = "Junior"
$ws.Range("G240").Value2 = "Deportes Tolima"
$ws.Range("H240").Value2 = 4
$ws.Range("I240").Value2 = 2
$ws.Range("J240").Value2 = "H"
$ws.Range("K240").Value2 = 1.95
$ws.Range("L240").Value2 = 3.3
$ws.Range("M240").Value2 = 4
$ws.Range("N240").Value2 = 1.909
$ws.Range("O240").Value2 = 3.75
$ws.Range("P240").Value2 = 3.8
$ws.Range("Q240").Value2 = -0.5
$ws.Range("R240").Value2 = 1.9
$ws.Range("S240").Value2 = 1.9
$ws.Range("T240").Value2 = 2.5
$ws.Range("U240").Value2 = 1.85
$ws.Range("V240").Value2 = 1.95
$ws.Range("W240").Value2 = 0.909
$ws.Range("X240").Value2 = -1
$ws.Range("Y240").Value2 = -1
$ws.Range("Z240").Value2 = 0.8999999999999999
$ws.Range("AA240").Value2 = -1
$ws.Range("AB240").Value2 = 0.8500000000000001
$ws.Range("AC240").Value2 = -1

# row 241
$ws.Range("B241").Value2 = 7528135
$ws.Range("F241").Value2 = "Independiente Medellin"
$ws.Range("G241").Value2 = "America de Cali"
$ws.Range("H241").Value2 = 2
$ws.Range("I241").Value2 = 1
$ws.Range("J241").Value2 = "H"
$ws.Range("K241").Value2 = 2.15
$ws.Range("L241").Value2 = 3.3
$ws.Range("M241").Value2 = 3.4
$ws.Range("N241").Value2 = 2.375
$ws.Range("O241").Value2 = 3.3
$ws.Range("P241").Value2 = 3.1
$ws.Range("Q241").Value2 = -0.25
$ws.Range("R241").Value2 = 2
$ws.Range("S241").Value2 = 1.8
$ws.Range("T241").Value2 = 2.5
$ws.Range("U241").Value2 = 1.975
$ws.Range("V241").Value2 = 1.825
$ws.Range("W241").Value2 = 1.375
$ws.Range("X241").Value2 = -1
$ws.Range("Y241").Value2 = -1
$ws.Range("Z241").Value2 = 1
$ws.Range("AA241").Value2 = -1
$ws.Range("AB241").Value2 = 0.9750000000000001
$ws.Range("AC241").Value2 = -1

# --- Existing rows 394-396 updated in place ---
# row 394
$ws.Range("B394").Value2 = 7658906
$ws.Range("C394").Value2 = "Colombia Primera A"
$ws.Range("D394").Value2 = "Colombia Primera A"
$ws.Range("E394").Value2 = 45395.66666666666
$ws.Range("F394").Value2 = "Fortaleza"
$ws.Range("G394").Value2 = "Once Caldas"
$ws.Range("H394").Value2 = 2
$ws.Range("I394").Value2 = 0
$ws.Range("J394").Value2 = "H"
$ws.Range("K394").Value2 = 2.1
$ws.Range("L394").Value2 = 3.1
$ws.Range("M394").Value2 = 3.8
$ws.Range("N394").Value2 = 1.909
$ws.Range("O394").Value2 = 3.3
$ws.Range("P394").Value2 = 4.5
$ws.Range("Q394").Value2 = -0.5
$ws.Range("R394").Value2 = 1.9
$ws.Range("S394").Value2 = 1.9
$ws.Range("T394").Value2 = 2
$ws.Range("U394").Value2 = 1.75
$ws.Range("V394").Value2 = 2.05
$ws.Range("W394").Value2 = 0.909
$ws.Range("X394").Value2 = -1
$ws.Range("Y394").Value2 = -1
$ws.Range("Z394").Value2 = 0.8999999999999999
$ws.Range("AA394").Value2 = -1
$ws.Range("AB394").Value2 = 0
$ws.Range("AC394").Value2 = -0

# row 395
$ws.Range("B395").Value2 = 7658964
$ws.Range("C395").Value2 = "Colombia Primera A"
$ws.Range("D395").Value2 = "Colombia Primera A"
$ws.Range("E395").Value2 = 45395.75694444445
$ws.Range("F395").Value2 = "Independiente Santa Fe"
$ws.Range("G395").Value2 = "Atletico Nacional Medellin"
$ws.Range("H395").Value2 = 0
$ws.Range("I395").Value2 = 0
$ws.Range("J395").Value2 = "D"
$ws.Range("K395").Value2 = 2.2
$ws.Range("L395").Value2 = 3.1
$ws.Range("M395").Value2 = 3.5
$ws.Range("N395").Value2 = 2
$ws.Range("O395").Value2 = 3.2
$ws.Range("P395").Value2 = 4.2
$ws.Range("Q395").Value2 = -0.5
$ws.Range("R395").Value2 = 2
$ws.Range("S395").Value2 = 1.8
$ws.Range("T395").Value2 = 2
$ws.Range("U395").Value2 = 1.85
$ws.Range("V395").Value2 = 1.95
$ws.Range("W395").Value2 = -1
$ws.Range("X395").Value2 = 2.2
$ws.Range("Y395").Value2 = -1
$ws.Range("Z395").Value2 = -1
$ws.Range("AA395").Value2 = 0.8
$ws.Range("AB395").Value2 = -1
$ws.Range("AC395").Value2 = 0.95

# row 396
$ws.Range("B396").Value2 = 7658966
$ws.Range("C396").Value2 = "Colombia Primera A"
$ws.Range("D396").Value2 = "Colombia Primera A"
$ws.Range("E396").Value2 = 45395.84722222222
$ws.Range("F396").Value2 = "Junior"
$ws.Range("G396").Value2 = "Envigado FC"
$ws.Range("H396").Value2 = 1
$ws.Range("I396").Value2 = 1
$ws.Range("J396").Value2 = "D"
$ws.Range("K396").Value2 = 1.4
$ws.Range("L396").Value2 = 4.75
$ws.Range("M396").Value2 = 7
$ws.Range("N396").Value2 = 1.333
$ws.Range("O396").Value2 = 5.25
$ws.Range("P396").Value2 = 10
$ws.Range("Q396").Value2 = -1.5
$ws.Range("R396").Value2 = 1.9
$ws.Range("S396").Value2 = 1.9
$ws.Range("T396").Value2 = 2.75
$ws.Range("U396").Value2 = 1.95
$ws.Range("V396").Value2 = 1.85
$ws.Range("W396").Value2 = -1
$ws.Range("X396").Value2 = 4.25
$ws.Range("Y396").Value2 = -1
$ws.Range("Z396").Value2 = -1
$ws.Range("AA396").Value2 = 0.8999999999999999
$ws.Range("AB396").Value2 = -1
$ws.Range("AC396").Value2 = 0.8500000000000001

# --- New rows 397-407 ---
# row 397
$ws.Range("A396").Copy() | Out-Null
$ws.Range("A397").PasteSpecial(-4122) | Out-Null
$ws.Range("E396").Copy() | Out-Null
$ws.Range("E397").PasteSpecial(-4122) | Out-Null
$ws.Range("B397").Value2 = 7658965
$ws.Range("C397").Value2 = "Colombia Primera A"
$ws.Range("D397").Value2 = "Colombia Primera A"
$ws.Range("E397").Value2 = 45395.9375
$ws.Range("F397").Value2 = "Deportivo Pereira"
$ws.Range("G397").Value2 = "Jaguares de Cordoba"
$ws.Range("H397").Value2 = 0
$ws.Range("I397").Value2 = 0
$ws.Range("J397").Value2 = "D"
$ws.Range("K397").Value2 = 1.615
$ws.Range("L397").Value2 = 3.6
$ws.Range("M397").Value2 = 5.75
$ws.Range("N397").Value2 = 1.65
$ws.Range("O397").Value2 = 3.5
$ws.Range("P397").Value2 = 6.5
$ws.Range("Q397").Value2 = -0.75
$ws.Range("R397").Value2 = 1.875
$ws.Range("S397").Value2 = 1.975
$ws.Range("T397").Value2 = 2.25
$ws.Range("U397").Value2 = 1.975
$ws.Range("V397").Value2 = 1.825
$ws.Range("W397").Value2 = -1
$ws.Range("X397").Value2 = 2.5
$ws.Range("Y397").Value2 = -1
$ws.Range("Z397").Value2 = -1
$ws.Range("AA397").Value2 = 0.9750000000000001
$ws.Range("AB397").Value2 = -1
$ws.Range("AC397").Value2 = 0.825

# row 398
$ws.Range("A397").Copy() | Out-Null
$ws.Range("A398").PasteSpecial(-4122) | Out-Null
$ws.Range("E397").Copy() | Out-Null
$ws.Range("E398").PasteSpecial(-4122) | Out-Null
$ws.Range("B398").Value2 = 7658977
$ws.Range("C398").Value2 = "Colombia Primera A"
$ws.Range("D398").Value2 = "Colombia Primera A"
$ws.Range("E398").Value2 = 45398.75
$ws.Range("F398").Value2 = "Envigado FC"
$ws.Range("G398").Value2 = "Deportes Tolima"
$ws.Range("K398").Value2 = 3.8
$ws.Range("L398").Value2 = 3.25
$ws.Range("M398").Value2 = 2.05
$ws.Range("N398").Value2 = 5.25
$ws.Range("O398").Value2 = 3.6
$ws.Range("P398").Value2 = 1.666
$ws.Range("Q398").Value2 = 0.75
$ws.Range("R398").Value2 = 1.9
$ws.Range("S398").Value2 = 1.95
$ws.Range("T398").Value2 = 2.25
$ws.Range("U398").Value2 = 1.875
$ws.Range("V398").Value2 = 1.975
$ws.Range("W398").Value2 = 0
$ws.Range("X398").Value2 = 0
$ws.Range("Y398").Value2 = 0
$ws.Range("Z398").Value2 = 0
$ws.Range("AA398").Value2 = 0

# row 399
$ws.Range("A398").Copy() | Out-Null
$ws.Range("A399").PasteSpecial(-4122) | Out-Null
$ws.Range("E398").Copy() | Out-Null
$ws.Range("E399").PasteSpecial(-4122) | Out-Null
$ws.Range("B399").Value2 = 7658973
$ws.Range("C399").Value2 = "Colombia Primera A"
$ws.Range("D399").Value2 = "Colombia Primera A"
$ws.Range("E399").Value2 = 45398.84027777778
$ws.Range("F399").Value2 = "Jaguares de Cordoba"
$ws.Range("G399").Value2 = "Deportivo Cali"
$ws.Range("K399").Value2 = 2.3
$ws.Range("L399").Value2 = 2.8
$ws.Range("M399").Value2 = 3.3
$ws.Range("N399").Value2 = 2.2
$ws.Range("O399").Value2 = 2.875
$ws.Range("P399").Value2 = 3.6
$ws.Range("Q399").Value2 = -0.25
$ws.Range("R399").Value2 = 1.875
$ws.Range("S399").Value2 = 1.975
$ws.Range("T399").Value2 = 2.25
$ws.Range("U399").Value2 = 2.05
$ws.Range("V399").Value2 = 1.8
$ws.Range("W399").Value2 = 0
$ws.Range("X399").Value2 = 0
$ws.Range("Y399").Value2 = 0
$ws.Range("Z399").Value2 = 0
$ws.Range("AA399").Value2 = 0

# row 400
$ws.Range("A399").Copy() | Out-Null
$ws.Range("A400").PasteSpecial(-4122) | Out-Null
$ws.Range("E399").Copy() | Out-Null
$ws.Range("E400").PasteSpecial(-4122) | Out-Null
$ws.Range("B400").Value2 = 7658907
$ws.Range("C400").Value2 = "Colombia Primera A"
$ws.Range("D400").Value2 = "Colombia Primera A"
$ws.Range("E400").Value2 = 45398.93055555555
$ws.Range("F400").Value2 = "Once Caldas"
$ws.Range("G400").Value2 = "Independiente Santa Fe"
$ws.Range("K400").Value2 = 2.5
$ws.Range("L400").Value2 = 2.875
$ws.Range("M400").Value2 = 2.9
$ws.Range("N400").Value2 = 3.4
$ws.Range("O400").Value2 = 2.875
$ws.Range("P400").Value2 = 2.25
$ws.Range("Q400").Value2 = 0.25
$ws.Range("R400").Value2 = 1.9
$ws.Range("S400").Value2 = 1.95
$ws.Range("T400").Value2 = 2
$ws.Range("U400").Value2 = 1.925
$ws.Range("V400").Value2 = 1.925
$ws.Range("W400").Value2 = 0
$ws.Range("X400").Value2 = 0
$ws.Range("Y400").Value2 = 0
$ws.Range("Z400").Value2 = 0
$ws.Range("AA400").Value2 = 0

# row 401
$ws.Range("A400").Copy() | Out-Null
$ws.Range("A401").PasteSpecial(-4122) | Out-Null
$ws.Range("E400").Copy() | Out-Null
$ws.Range("E401").PasteSpecial(-4122) | Out-Null
$ws.Range("B401").Value2 = 7736843
$ws.Range("C401").Value2 = "Colombia Primera A"
$ws.Range("D401").Value2 = "Colombia Primera A"
$ws.Range("E401").Value2 = 45399.75
$ws.Range("F401").Value2 = "Aguilas Doradas"
$ws.Range("G401").Value2 = "Alianza"
$ws.Range("K401").Value2 = 1.727
$ws.Range("L401").Value2 = 3.4
$ws.Range("M401").Value2 = 4.75
$ws.Range("N401").Value2 = 1.666
$ws.Range("O401").Value2 = 3.4
$ws.Range("P401").Value2 = 5
$ws.Range("Q401").Value2 = -0.75
$ws.Range("R401").Value2 = 2
$ws.Range("S401").Value2 = 1.85
$ws.Range("T401").Value2 = 2.25
$ws.Range("U401").Value2 = 2.05
$ws.Range("V401").Value2 = 1.8
$ws.Range("W401").Value2 = 0
$ws.Range("X401").Value2 = 0
$ws.Range("Y401").Value2 = 0
$ws.Range("Z401").Value2 = 0
$ws.Range("AA401").Value2 = 0

# row 402
$ws.Range("A401").Copy() | Out-Null
$ws.Range("A402").PasteSpecial(-4122) | Out-Null
$ws.Range("E401").Copy() | Out-Null
$ws.Range("E402").PasteSpecial(-4122) | Out-Null
$ws.Range("B402").Value2 = 7658971
$ws.Range("C402").Value2 = "Colombia Primera A"
$ws.Range("D402").Value2 = "Colombia Primera A"
$ws.Range("E402").Value2 = 45399.84027777778
$ws.Range("F402").Value2 = "Millonarios"
$ws.Range("G402").Value2 = "Junior"
$ws.Range("K402").Value2 = 2
$ws.Range("L402").Value2 = 3.2
$ws.Range("M402").Value2 = 3.75
$ws.Range("N402").Value2 = 2
$ws.Range("O402").Value2 = 3.2
$ws.Range("P402").Value2 = 3.75
$ws.Range("Q402").Value2 = -0.5
$ws.Range("R402").Value2 = 2.025
$ws.Range("S402").Value2 = 1.825
$ws.Range("T402").Value2 = 2.25
$ws.Range("U402").Value2 = 1.95
$ws.Range("V402").Value2 = 1.9
$ws.Range("W402").Value2 = 0
$ws.Range("X402").Value2 = 0
$ws.Range("Y402").Value2 = 0
$ws.Range("Z402").Value2 = 0
$ws.Range("AA402").Value2 = 0

# row 403
$ws.Range("A402").Copy() | Out-Null
$ws.Range("A403").PasteSpecial(-4122) | Out-Null
$ws.Range("E402").Copy() | Out-Null
$ws.Range("E403").PasteSpecial(-4122) | Out-Null
$ws.Range("B403").Value2 = 7658975
$ws.Range("C403").Value2 = "Colombia Primera A"
$ws.Range("D403").Value2 = "Colombia Primera A"
$ws.Range("E403").Value2 = 45399.93055555555
$ws.Range("F403").Value2 = "Atletico Nacional Medellin"
$ws.Range("G403").Value2 = "Deportivo Pereira"
$ws.Range("K403").Value2 = 1.75
$ws.Range("L403").Value2 = 3.5
$ws.Range("M403").Value2 = 4.5
$ws.Range("N403").Value2 = 1.75
$ws.Range("O403").Value2 = 3.5
$ws.Range("P403").Value2 = 4.5
$ws.Range("Q403").Value2 = -0.75
$ws.Range("R403").Value2 = 2.05
$ws.Range("S403").Value2 = 1.8
$ws.Range("T403").Value2 = 2.25
$ws.Range("U403").Value2 = 1.975
$ws.Range("V403").Value2 = 1.875
$ws.Range("W403").Value2 = 0
$ws.Range("X403").Value2 = 0
$ws.Range("Y403").Value2 = 0
$ws.Range("Z403").Value2 = 0
$ws.Range("AA403").Value2 = 0

# row 404
$ws.Range("A403").Copy() | Out-Null
$ws.Range("A404").PasteSpecial(-4122) | Out-Null
$ws.Range("E403").Copy() | Out-Null
$ws.Range("E404").PasteSpecial(-4122) | Out-Null
$ws.Range("B404").Value2 = 7658908
$ws.Range("C404").Value2 = "Colombia Primera A"
$ws.Range("D404").Value2 = "Colombia Primera A"
$ws.Range("E404").Value2 = 45400.75
$ws.Range("F404").Value2 = "La Equidad"
$ws.Range("G404").Value2 = "Fortaleza"
$ws.Range("K404").Value2 = 1.833
$ws.Range("L404").Value2 = 3.25
$ws.Range("M404").Value2 = 4.333
$ws.Range("N404").Value2 = 1.833
$ws.Range("O404").Value2 = 3.25
$ws.Range("P404").Value2 = 4.333
$ws.Range("Q404").Value2 = -0.5
$ws.Range("R404").Value2 = 1.85
$ws.Range("S404").Value2 = 2
$ws.Range("T404").Value2 = 2.25
$ws.Range("U404").Value2 = 2
$ws.Range("V404").Value2 = 1.85
$ws.Range("W404").Value2 = 0
$ws.Range("X404").Value2 = 0
$ws.Range("Y404").Value2 = 0
$ws.Range("Z404").Value2 = 0
$ws.Range("AA404").Value2 = 0

# row 405
$ws.Range("A404").Copy() | Out-Null
$ws.Range("A405").PasteSpecial(-4122) | Out-Null
$ws.Range("E404").Copy() | Out-Null
$ws.Range("E405").PasteSpecial(-4122) | Out-Null
$ws.Range("B405").Value2 = 7658972
$ws.Range("C405").Value2 = "Colombia Primera A"
$ws.Range("D405").Value2 = "Colombia Primera A"
$ws.Range("E405").Value2 = 45400.84027777778
$ws.Range("F405").Value2 = "Patriotas FC"
$ws.Range("G405").Value2 = "Independiente Medellin"
$ws.Range("K405").Value2 = 2.625
$ws.Range("L405").Value2 = 3.1
$ws.Range("M405").Value2 = 2.625
$ws.Range("N405").Value2 = 2.625
$ws.Range("O405").Value2 = 3.1
$ws.Range("P405").Value2 = 2.625
$ws.Range("Q405").Value2 = 0
$ws.Range("R405").Value2 = 1.875
$ws.Range("S405").Value2 = 1.975
$ws.Range("T405").Value2 = 2.25
$ws.Range("U405").Value2 = 2.05
$ws.Range("V405").Value2 = 1.8
$ws.Range("W405").Value2 = 0
$ws.Range("X405").Value2 = 0
$ws.Range("Y405").Value2 = 0
$ws.Range("Z405").Value2 = 0
$ws.Range("AA405").Value2 = 0

# row 406
$ws.Range("A405").Copy() | Out-Null
$ws.Range("A406").PasteSpecial(-4122) | Out-Null
$ws.Range("E405").Copy() | Out-Null
$ws.Range("E406").PasteSpecial(-4122) | Out-Null
$ws.Range("B406").Value2 = 7658976
$ws.Range("C406").Value2 = "Colombia Primera A"
$ws.Range("D406").Value2 = "Colombia Primera A"
$ws.Range("E406").Value2 = 45400.93055555555
$ws.Range("F406").Value2 = "America de Cali"
$ws.Range("G406").Value2 = "Deportivo Pasto"
$ws.Range("K406").Value2 = 1.533
$ws.Range("L406").Value2 = 3.6
$ws.Range("M406").Value2 = 6.5
$ws.Range("N406").Value2 = 1.5
$ws.Range("O406").Value2 = 3.75
$ws.Range("P406").Value2 = 7
$ws.Range("Q406").Value2 = -1
$ws.Range("R406").Value2 = 1.9
$ws.Range("S406").Value2 = 1.95
$ws.Range("T406").Value2 = 2.25
$ws.Range("U406").Value2 = 1.85
$ws.Range("V406").Value2 = 2
$ws.Range("W406").Value2 = 0
$ws.Range("X406").Value2 = 0
$ws.Range("Y406").Value2 = 0
$ws.Range("Z406").Value2 = 0
$ws.Range("AA406").Value2 = 0

# row 407
$ws.Range("A406").Copy() | Out-Null
$ws.Range("A407").PasteSpecial(-4122) | Out-Null
$ws.Range("E406").Copy() | Out-Null
$ws.Range("E407").PasteSpecial(-4122) | Out-Null
$ws.Range("B407").Value2 = 7658974
$ws.Range("C407").Value2 = "Colombia Primera A"
$ws.Range("D407").Value2 = "Colombia Primera A"
$ws.Range("E407").Value2 = 45401.92361111111
$ws.Range("F407").Value2 = "Atletico Bucaramanga"
$ws.Range("G407").Value2 = "Boyaca Chico"
$ws.Range("K407").Value2 = 1.65
$ws.Range("L407").Value2 = 3.6
$ws.Range("M407").Value2 = 5
$ws.Range("N407").Value2 = 1.666
$ws.Range("O407").Value2 = 3.6
$ws.Range("P407").Value2 = 5
$ws.Range("Q407").Value2 = -0.75
$ws.Range("R407").Value2 = 1.875
$ws.Range("S407").Value2 = 1.975
$ws.Range("T407").Value2 = 2.25
$ws.Range("U407").Value2 = 1.925
$ws.Range("V407").Value2 = 1.925
$ws.Range("W407").Value2 = 0
$ws.Range("X407").Value2 = 0
$ws.Range("Y407").Value2 = 0
$ws.Range("Z407").Value2 = 0
$ws.Range("AA407").Value2 = 0

